$d = $word.ActiveDocument

# The target paragraph that currently ends the body text ("- Tipo atividade...")
$anchor = $d.Paragraphs.Item(22)
$cursor = $anchor.Range
$cursor.Collapse(0)

# Insert three blank paragraphs
$cursor.InsertParagraphAfter()
$cursor.Collapse(0)
$cursor.InsertParagraphAfter()
$cursor.Collapse(0)
$cursor.InsertParagraphAfter()
$cursor.Collapse(0)

# Insert the paragraph that will hold:
#   "- " + "Parte para cadastro de empresas que fazem doaç" + "ão" (special formatting)
$cursor.InsertParagraphAfter()
$cursor.Collapse(0)

$pDoacao = $d.Paragraphs.Item(26)
$rDoacao = $pDoacao.Range
$rDoacao.InsertAfter("- ")
$rDoacao.InsertAfter("Parte para cadastro de empresas que fazem doaç")
$splitPos = $rDoacao.End - 1
$rDoacao.InsertAfter("ão")

$rAo = $d.Range($splitPos, $rDoacao.End - 1)
$rAo.Font.NameFarEast = "Noto Serif CJK SC"
$rAo.Font.NameBi = "Lohit Devanagari"
$rAo.Font.Color = -16777216
$rAo.Font.Kerning = 1
$rAo.Font.Size = 12
$rAo.Font.SizeBi = 12
$rAo.LanguageID = "pt-BR"
$rAo.LanguageIDFarEast = "zh-CN"
$rAo.LanguageIDOther = "hi-IN"

# Next paragraph: "- parte para representante legal da empresa" (same special formatting)
$pRep = $d.Paragraphs.Item(27)
$rRep = $pRep.Range
$rRep.InsertAfter("- parte para representante legal da empresa")
$rRepFmt = $d.Range($rRep.Start, $rRep.End - 1)
$rRepFmt.Font.NameFarEast = "Noto Serif CJK SC"
$rRepFmt.Font.NameBi = "Lohit Devanagari"
$rRepFmt.Font.Color = -16777216
$rRepFmt.Font.Kerning = 1
$rRepFmt.Font.Size = 12
$rRepFmt.Font.SizeBi = 12
$rRepFmt.LanguageID = "pt-BR"
$rRepFmt.LanguageIDFarEast = "zh-CN"
$rRepFmt.LanguageIDOther = "hi-IN"
